$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 216, pushing existing rows 216:263 down to 217:264
$ws.Rows("216:216").Insert()

# Populate the newly inserted row 216 with the new weekly price record
$ws.Range("A216").Value = 4
$ws.Range("B216").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C216").Value = "Los Lagos"
$ws.Range("D216").Value = 44711
$ws.Range("E216").Value = 10
$ws.Range("F216").Value = 100112021
$ws.Range("G216").Value = "Ají"
$ws.Range("H216").Value = "Inferno"
$ws.Range("I216").Value = "Primera"
$ws.Range("J216").Value = 35
$ws.Range("K216").Value = 30000
$ws.Range("L216").Value = 30000
$ws.Range("M216").Value = 30000
$ws.Range("N216").Value = "$/caja 12 kilos"
$ws.Range("O216").Value = "Región de Arica y Parinacota"
$ws.Range("P216").Value = 2500
$ws.Range("Q216").Value = 12
$ws.Range("R216").Value = "Hortaliza"
